# Apply hybrid bold + color (2C3E50) highlighting to quantitative impact
# metrics (percentages, dollar amounts, large numbers) across the resume.
#
# Strategy: for each target paragraph, scope a Find/Replace search to that
# paragraph's Range so we don't accidentally match the same numeric string
# appearing elsewhere in the document (e.g. "23% to 64%" also shows up in
# the Professional Summary and Key Projects sections, which must stay
# untouched). Each match is located with Range.Find.Execute, which collapses
# the Range to the found text, and then Font.Bold / Font.Color are applied
# directly to that sub-range -- Word automatically splits the run and
# preserves the xml:space="preserve" attribute where needed.

$HighlightColor = 5258796   # BGR-encoded value of RGB 2C3E50 (w:val="2C3E50")

function Apply-MetricHighlight {
    param(
        $Paragraph,
        [string]$SearchText
    )

    $rng = $Paragraph.Range
    $found = $rng.Find.Execute($SearchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = $true
        $rng.Font.Color = $HighlightColor
    }
    Write-Host "Highlighted '$SearchText':" $found
}

$d = $word.ActiveDocument

# • Discovered systematic race coding errors ... accuracy from 23% to 64%
Apply-MetricHighlight $d.Paragraphs.Item(10) "23%"
Apply-MetricHighlight $d.Paragraphs.Item(10) "64%"

# • Utilized advanced sampling methods ... from ±4.2% to ±2.1% ... from 71% to 87% ...
Apply-MetricHighlight $d.Paragraphs.Item(12) "±4.2%"
Apply-MetricHighlight $d.Paragraphs.Item(12) "±2.1%"
Apply-MetricHighlight $d.Paragraphs.Item(12) "71%"
Apply-MetricHighlight $d.Paragraphs.Item(12) "87%"

# • Trigonometric algorithm ... reduced mapping costs by 73.5%, saving ... $4.7M ...
Apply-MetricHighlight $d.Paragraphs.Item(13) "73.5%"
Apply-MetricHighlight $d.Paragraphs.Item(13) "$4.7M"

# • Built real-time FEC analysis systems ... valued over $2 trillion
Apply-MetricHighlight $d.Paragraphs.Item(14) "$2"

# • Modernized legacy ETL processes ... reducing processing time by 57%
Apply-MetricHighlight $d.Paragraphs.Item(24) "57%"

# • Revenue generation: Delivered $4.9M additional revenue through optimization
Apply-MetricHighlight $d.Paragraphs.Item(50) "$4.9M"

# • 23% conversion rate improvement
Apply-MetricHighlight $d.Paragraphs.Item(51) "23%"

# • Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations
Apply-MetricHighlight $d.Paragraphs.Item(53) "12,847"
